$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (row 1)
$ws.Range("A1").Value = "angular frequency"
$ws.Range("B1").Value = "Z_imag"
$ws.Range("C1").Value = "Z_real"

# Update the selected cell to match the saved view state
$ws.Range("C3").Select()
